$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Admin_SZ_Portuguese")

# --- Row 6: finish out the "Room Tax" row and add the Room Types columns ---
# Reuse existing shared strings first (order doesn't affect the shared-string table).
$ws.Cells.Item(6, 5).Value  = "Message_Text1"          # E6
$ws.Cells.Item(6, 7).Value  = "Message_Text2"          # G6
$ws.Cells.Item(6, 16).Value = " salvo com sucesso"     # P6

# New strings - the order below reproduces the workbook's shared-string insertion order.
$ws.Cells.Item(6, 8).Value  = " salva com sucesso"     # H6
$ws.Cells.Item(6, 9).Value  = "Status_src"             # I6
$ws.Cells.Item(6, 10).Value = "on.GIF"                 # J6
$ws.Cells.Item(6, 11).Value = "RoomTypesList_Title"    # K6
$ws.Cells.Item(6, 13).Value = "Message_Text3"          # M6
$ws.Cells.Item(6, 15).Value = "Message_Text4"          # O6
$ws.Cells.Item(6, 6).Value  = "'- Taxa de quarto "      # F6 (leading ' => quote-prefix style, like F3:F5)
$ws.Cells.Item(6, 12).Value = "Lista de tipos de quarto"   # L6
$ws.Cells.Item(6, 14).Value = "'- Tipo de quarto (s) "  # N6 (leading ' => quote-prefix style)

# --- Row 7: AddOns list test case ---
$ws.Cells.Item(7, 1).Value = "fn_verifyAdditionOfAddOns"   # A7
$ws.Cells.Item(7, 3).Value = "AddOnsList_Title"             # C7
$ws.Cells.Item(7, 4).Value = "Lista de adições"             # D7
$ws.Cells.Item(7, 5).Value = "Status_src"                   # E7
$ws.Cells.Item(7, 6).Value = "'on.GIF"                       # F7 (quote-prefix style)

# --- Row 8: Add package with inclusion test case ---
$ws.Cells.Item(8, 1).Value = "fn_verifyAddAPackageWithInclusion"  # A8
$ws.Cells.Item(8, 3).Value = "AddInclusion_Title"                  # C8
$ws.Cells.Item(8, 4).Value = "Adicionar inclusões"                 # D8
$ws.Cells.Item(8, 5).Value = "Status_src"                          # E8
$ws.Cells.Item(8, 6).Value = "'on.GIF"                              # F8 (quote-prefix style)

# --- Column widths for the new K:P columns ---
$ws.Columns.Item(11).ColumnWidth = 18.666666666666668   # K -> 19.5703125
$ws.Columns.Item(12).ColumnWidth = 22.0                 # L -> 22.85546875
$ws.Columns.Item(13).ColumnWidth = 13.833333333333334   # M -> 14.7109375
$ws.Columns.Item(14).ColumnWidth = 17.666666666666668   # N -> 18.42578125
$ws.Columns.Item(15).ColumnWidth = 13.833333333333334   # O -> 14.7109375
$ws.Columns.Item(16).ColumnWidth = 16.833333333333332   # P -> 17.7109375

# --- Selection moves to D21 ---
$ws.Range("D21").Select()
